$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Helper: force a Range to remain its own <w:r> run (instead of being
# re-merged into a neighbouring run with identical formatting) by
# toggling a character property on and back off again.
# ------------------------------------------------------------------
function Lock-Run($rng) {
    $rng.Bold = 1
    $rng.Bold = 0
}

# ====================================================================
# Part 1: "Výsledné video nastavte jako " -> "Výsledné videa nastavte jako "
#   becomes three runs: "Vysledne vide" | "a" | " nastavte jako "
# ====================================================================
$find1 = $d.Content
$find1.Find.Execute("Výsledné video nastavte jako")
$para1Start = $find1.Start

# The "o" that must turn into "a" is the 14th character (0-based index 13)
# of "Výsledné video nastavte jako" ("Výsledné vide" + "o" + " nastavte jako").
$oStart = $para1Start + 13
$oRange = $d.Range($oStart, $oStart + 1)
$oRange.Text = "a"

# Re-fetch the (now 1-char) range for "a" and lock it into its own run,
# separate from "Vysledne vide" before it and " nastavte jako " after it.
$aRange = $d.Range($oStart, $oStart + 1)
Lock-Run $aRange

# ====================================================================
# Part 2: " a odevzdejte jako odkaz do " -> " a odevzdejte odkazy na videa do "
#   becomes five runs: " a odevzdejte " | "od" | "kaz" | "y na videa" | " do "
# ====================================================================
$find2 = $d.Content
$find2.Find.Execute(" a odevzdejte jako odkaz do ")
$seg2Start = $find2.Start

# " a odevzdejte jako odkaz do "
#  0123456789...
# Remove "jako " (characters 14-18 inclusive, i.e. offset 14..19)
$jakoStart = $seg2Start + 14
$jakoEnd = $jakoStart + 5
$jakoRange = $d.Range($jakoStart, $jakoEnd)
$jakoRange.Text = ""

# After removal: " a odevzdejte odkaz do ", "odkaz" starts right where
# "jako " used to start.
$odkazStart = $jakoStart

# Insert "y na videa" right after "odkaz" (still a single run at this
# point), turning it into "odkazy na videa".
$insertPoint = $d.Range($odkazStart + 5, $odkazStart + 5)
$insertPoint.InsertAfter("y na videa")

# InsertAfter re-normalises/merges the surrounding text back into one
# run, so now re-split it into "od" | "kaz" | "y na videa" | " do ".
$odEnd = $odkazStart + 2
$kazEnd = $odEnd + 3
$yEnd = $kazEnd + 10

$odRange = $d.Range($odkazStart, $odEnd)
Lock-Run $odRange

$kazRange = $d.Range($odEnd, $kazEnd)
Lock-Run $kazRange

$yRange = $d.Range($kazEnd, $yEnd)
Lock-Run $yRange
